# Update Zeit + Statusprotokoll
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add new rows of task descriptions in column B for rows 11-14
# Set in this order so shared-string table indices come out as: 9,10,11,12
$ws.Range("B11").Value = "Recherche für Erweiterungen"
$ws.Range("B13").Value = "Implementierung von Erweiterungen im Frontend"
$ws.Range("B14").Value = "Änderung der Datenstruktur im Frontend"
$ws.Range("B12").Value = "Fixing der Forms für das Mitarbeiterupdate"

# Update the active selection to B12
$ws.Range("B12").Select()
